# fix bug api with time system
# Remove the two obsolete course rows "Selected Topics in Computer Engineering
# and Informatics" (03603496) and "Special Problems" (03603498) from the
# course list on Sheet1. Deleting the rows shifts the subsequent rows
# (Co-operative Education Preparation / Co-operative Education) up, which
# also drops the now-unused shared strings automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete rows 65 and 66 (03603496 / 03603498), shifting rows below them up.
$ws.Range("A65:D66").EntireRow.Delete()
